# Added user option for defining continuity
# Adds two new columns (start_continuity / end_continuity) to the
# Trend_instructions sheet, each populated with "F" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

# New headers in L1:M1
$ws.Range("L1").Value = "start_continuity"
$ws.Range("M1").Value = "end_continuity"

# New values in L2:M5 (one per existing data row)
$ws.Range("L2:L5").Value = "F"
$ws.Range("M2:M5").Value = "F"

# Match the bestFit-style width used by the sheet's other text columns
$ws.Columns.Item(12).ColumnWidth = 13.1640625

# Update the view/selection state to reflect scrolling to the new columns
$ws.Activate()
$ws.Range("C1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("L8").Select()
